$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = '@'
$c.Value = '29.725.62'
$c.Style = 'Normal'
$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.02%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = '@'
$c.Value = '1.885.17'
$c.Style = 'Normal'

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.Style = 'Normal'
$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.20%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '0.7921'
$c.Style = 'Normal'
$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = '@'
$c.Value = '  -6.41%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = '@'
$c.Value = '241.13'
$c.Style = 'Normal'
$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.40%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = '@'
$c.Value = '1.0000'
$c.Style = 'Normal'
$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.09%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = '@'
$c.Value = '0.3163'
$c.Style = 'Normal'
$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = '@'
$c.Value = '  -2.72%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = '@'
$c.Value = '25.35'
$c.Style = 'Normal'
$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = '@'
$c.Value = '  -5.45%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = '@'
$c.Value = '0.06985'
$c.Style = 'Normal'
$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.23%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.06%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = '@'
$c.Value = '0.7607'
$c.Style = 'Normal'
$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.84%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = '@'
$c.Value = '1.896.61'
$c.Style = 'Normal'
$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.58%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = '@'
$c.Value = '5.292'
$c.Style = 'Normal'
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = '@'
$c.Value = '  +1.17%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Value = '92.09'
$c.Style = 'Normal'
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.98%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = '@'
$c.Value = '29.704.67'
$c.Style = 'Normal'
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.09%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = '@'
$c.Value = '13.82'
$c.Style = 'Normal'
$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = '@'
$c.Value = '  -2.65%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = '@'
$c.Value = '5.933'
$c.Style = 'Normal'
$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.65%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = '@'
$c.Value = '243.00'
$c.Style = 'Normal'
$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.85%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = '@'
$c.Value = '0.000007666'
$c.Style = 'Normal'
$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.52%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '8.177'
$c.Style = 'Normal'
$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = '@'
$c.Value = '  +16.42%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = '@'
$c.Value = '0.9994'
$c.Style = 'Normal'

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = '2.133.34'
$c.Style = 'Normal'
$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.92%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = '@'
$c.Value = '1.002'
$c.Style = 'Normal'
$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.30%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = '@'
$c.Value = '0.1674'
$c.Style = 'Normal'
$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = '@'
$c.Value = '  +3.03%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Value = '9.274'
$c.Style = 'Normal'
$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.09%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Value = '164.17'
$c.Style = 'Normal'
$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = '@'
$c.Value = '  -3.19%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = '@'
$c.Value = '18.59'
$c.Style = 'Normal'
$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = '@'
$c.Value = '  -2.03%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = '@'
$c.Value = '2.045'
$c.Style = 'Normal'
$c = $ws.Cells.Item(29, 5)
$c.NumberFormat = '@'
$c.Value = '  -2.16%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(30, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.87%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = '@'
$c.Value = '1.531'
$c.Style = 'Normal'
$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = '@'
$c.Value = '  +1.03%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = '@'
$c.Value = '4.367'
$c.Style = 'Normal'
$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = '@'
$c.Value = '  +1.28%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '0.05661'
$c.Style = 'Normal'
$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.32%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = '@'
$c.Value = '4.045'
$c.Style = 'Normal'
$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.32%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = '@'
$c.Value = '1.258'
$c.Style = 'Normal'
$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = '@'
$c.Value = '  -2.46%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = '@'
$c.Value = '0.7328'
$c.Style = 'Normal'
$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.67%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = '@'
$c.Value = '0.9956'
$c.Style = 'Normal'
$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.26%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = '@'
$c.Value = '  -3.87%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = '@'
$c.Value = '0.01906'
$c.Style = 'Normal'
$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.65%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.87%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = '@'
$c.Value = '0.4394'
$c.Style = 'Normal'
$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.21%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = '@'
$c.Value = '72.25'
$c.Style = 'Normal'
$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.46%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Value = '5.810'
$c.Style = 'Normal'
$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = '@'
$c.Value = '  -3.49%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = '@'
$c.Value = '0.9992'
$c.Style = 'Normal'
$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.01%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value = '0.8349'
$c.Style = 'Normal'

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = '@'
$c.Value = '102.42'
$c.Style = 'Normal'
$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = '@'
$c.Value = '  +1.00%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = '@'
$c.Value = '1.020.44'
$c.Style = 'Normal'
$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = '@'
$c.Value = '  +3.10%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(48, 2)
$c.NumberFormat = '@'
$c.Value = 'EnergySwap'
$c.Style = 'Normal'
$c = $ws.Cells.Item(48, 3)
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c.Style = 'Normal'
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = '@'
$c.Value = '9.924'
$c.Style = 'Normal'
$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = '@'
$c.Value = '  +1.44%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(49, 2)
$c.NumberFormat = '@'
$c.Value = 'RenderToken'
$c.Style = 'Normal'
$c = $ws.Cells.Item(49, 3)
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c.Style = 'Normal'
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = '@'
$c.Value = '1.859'
$c.Style = 'Normal'
$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = '@'
$c.Value = '  -2.63%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Value = '7.404'
$c.Style = 'Normal'
$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = '@'
$c.Value = '  -2.94%  '
$c.Style = 'Normal'

$c = $ws.Cells.Item(51, 2)
$c.NumberFormat = '@'
$c.Value = 'RocketPoolETH'
$c.Style = 'Normal'
$c = $ws.Cells.Item(51, 3)
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$c.Style = 'Normal'
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = '@'
$c.Value = '2.007.69'
$c.Style = 'Normal'
$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = '@'
$c.Value = '  -2.59%  '
$c.Style = 'Normal'
